# Scheduled market-data refresh: update per-job Leve profit sheets with
# freshly pulled currentAveragePrice* figures (and downstream Leve*
# price/profit columns), and drop stale zero/placeholder rows on BSM
# whose source items no longer resolved on this run.

$wb = $excel.ActiveWorkbook

# --- ALC ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 418.07693
$ws.Range("I107").Value = 322.5
$ws.Range("J107").Value = 736.6667
$ws.Range("K107").Value = 322.5
$ws.Range("L107").Value = 736.6667
$ws.Range("M107").Value = 1597.5
$ws.Range("N107").Value = -4576.6667

# --- ARM -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2619.2693
$ws.Range("I74").Value = 705
$ws.Range("K74").Value = 705
$ws.Range("M74").Value = 169

$ws.Range("H77").Value = 2619.2693
$ws.Range("I77").Value = 705
$ws.Range("K77").Value = 3525
$ws.Range("M77").Value = 843

# --- BSM: clear stale H:N figures for rows whose leve data dropped out ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H117:N120").ClearContents()
$ws.Range("H122:N135").ClearContents()
$ws.Range("H137:N141").ClearContents()

# --- CRP -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1114.5454
$ws.Range("I58").Value = 995.7143
$ws.Range("K58").Value = 995.7143
$ws.Range("M58").Value = -792.7143

$ws.Range("H136").Value = 1114.5454
$ws.Range("I136").Value = 995.7143
$ws.Range("K136").Value = 2987.1429
$ws.Range("M136").Value = -437.1428999999998

# --- CUL -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1081.1082
$ws.Range("I122").Value = 570.5
$ws.Range("J122").Value = 1179.9354
$ws.Range("K122").Value = 5134.5
$ws.Range("L122").Value = 10619.4186
$ws.Range("M122").Value = -2684.5
$ws.Range("N122").Value = -15519.4186

# --- GSM -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 6302.375
$ws.Range("I43").Value = 1625
$ws.Range("J43").Value = 10979.75
$ws.Range("K43").Value = 1625
$ws.Range("L43").Value = 10979.75
$ws.Range("M43").Value = -1474
$ws.Range("N43").Value = -11281.75

$ws.Range("H113").Value = 4933.84
$ws.Range("I113").Value = 6964.75
$ws.Range("J113").Value = 1323.3334
$ws.Range("K113").Value = 6964.75
$ws.Range("L113").Value = 1323.3334
$ws.Range("M113").Value = -4794.75
$ws.Range("N113").Value = -5663.3334

# --- LTW -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1863.4783
$ws.Range("I68").Value = 1787.5
$ws.Range("J68").Value = 2037.1428
$ws.Range("K68").Value = 1787.5
$ws.Range("L68").Value = 2037.1428
$ws.Range("M68").Value = -1038.5
$ws.Range("N68").Value = -3535.1428

$ws.Range("H71").Value = 1863.4783
$ws.Range("I71").Value = 1787.5
$ws.Range("J71").Value = 2037.1428
$ws.Range("K71").Value = 8937.5
$ws.Range("L71").Value = 10185.714
$ws.Range("M71").Value = -5193.5
$ws.Range("N71").Value = -17673.714

$ws.Range("H100").Value = 1027.5714
$ws.Range("I100").Value = 900.5454999999999
$ws.Range("J100").Value = 1493.3334
$ws.Range("K100").Value = 900.5454999999999
$ws.Range("L100").Value = 1493.3334
$ws.Range("M100").Value = -359.5454999999999
$ws.Range("N100").Value = -2575.3334

# --- WVR -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2929.4263
$ws.Range("I132").Value = 3037.9246
$ws.Range("J132").Value = 2210.625
$ws.Range("K132").Value = 9113.773799999999
$ws.Range("L132").Value = 6631.875
$ws.Range("M132").Value = -6583.773799999999
$ws.Range("N132").Value = -11691.875
